$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    11 = -12.379
    12 = -10.984
    15 = -13.599
    27 = -13.096
    28 = -13.012
    31 = -13.589
    32 = -13.895
    36 = -12.732
    38 = -12.703
    46 = -14.033
    54 = -13.408
    55 = -13.599
    56 = -13.544
    67 = -11.595
    69 = -11.038
    72 = -11.555
    73 = -12.515
    83 = -13.21
    86 = -13.665
    91 = -10.98
    93 = -11.979
    99 = -12.715
}

foreach ($row in $updates.Keys) {
    $ws.Range("C$row").Value = $updates[$row]
}
